$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sparse odds corrections across rows 3, 5, 6, 7, 8, 9, 10 ---
# Row 3
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
# Row 5
$ws.Range("G5").Value = 3.4
$ws.Range("I5").Value = 2.3
$ws.Range("Q5").Value = 2.7
$ws.Range("R5").Value = 1.44
$ws.Range("AH5").Value = 9.5
$ws.Range("AU5").Value = 9.5
$ws.Range("AX5").Value = 15
# Row 6
$ws.Range("N6").Value = 10
# Row 7
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("W7").Value = 6
$ws.Range("AM7").Value = 1250
$ws.Range("AR7").Value = 51
# Row 8
$ws.Range("G8").Value = 4.33
$ws.Range("I8").Value = 1.85
$ws.Range("L8").Value = 2.6
$ws.Range("U8").Value = 2.2
$ws.Range("V8").Value = 1.62
$ws.Range("Z8").Value = 51
$ws.Range("AE8").Value = 21
$ws.Range("AH8").Value = 7.5
$ws.Range("AR8").Value = 151
# Row 9
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("AM9").Value = 800
# Row 10
$ws.Range("G10").Value = 2.25
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 3.4
$ws.Range("Q10").Value = 2.3
$ws.Range("R10").Value = 1.6
$ws.Range("U10").Value = 1.91
$ws.Range("V10").Value = 1.8
$ws.Range("X10").Value = 10
$ws.Range("Z10").Value = 21
$ws.Range("AA10").Value = 21
$ws.Range("AE10").Value = 15
$ws.Range("AG10").Value = 9
$ws.Range("AH10").Value = 15
$ws.Range("AJ10").Value = 34
$ws.Range("AK10").Value = 29
$ws.Range("AU10").Value = 8

# --- Row 13: Portugal AVS vs FC Porto replaced by USA MLS FC Cincinnati vs New York City ---
$ws.Range("A13").Value = "GAw0YMbl"
$ws.Range("B13").Value = "28/10/2024"
$ws.Range("C13").Value = "19:45"
$ws.Range("D13").Value = "USA - MLS"
$ws.Range("E13").Value = "FC Cincinnati"
$ws.Range("F13").Value = "New York City"
$ws.Range("G13").Value = 1.91
$ws.Range("H13").Value = 3.75
$ws.Range("I13").Value = 3.75
$ws.Range("J13").Value = 2.4
$ws.Range("K13").Value = 2.4
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 15
$ws.Range("O13").Value = 1.18
$ws.Range("P13").Value = 4.5
$ws.Range("Q13").Value = 1.62
$ws.Range("R13").Value = 2.25
$ws.Range("S13").Value = 1.29
$ws.Range("T13").Value = 3.5
$ws.Range("U13").Value = 1.53
$ws.Range("V13").Value = 2.38
$ws.Range("W13").Value = 11
$ws.Range("X13").Value = 11
$ws.Range("Y13").Value = 9
$ws.Range("Z13").Value = 17
$ws.Range("AA13").Value = 13
$ws.Range("AB13").Value = 21
$ws.Range("AC13").Value = 17
$ws.Range("AD13").Value = 7.5
$ws.Range("AE13").Value = 12
$ws.Range("AF13").Value = 34
$ws.Range("AG13").Value = 15
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 13
$ws.Range("AJ13").Value = 41
$ws.Range("AK13").Value = 26
$ws.Range("AL13").Value = 29
$ws.Range("AM13").Value = 101
$ws.Range("AN13").Value = 4.33
$ws.Range("AO13").Value = 9.5
$ws.Range("AP13").Value = 17
$ws.Range("AQ13").Value = 29
$ws.Range("AR13").Value = 41
$ws.Range("AS13").Value = 81
$ws.Range("AT13").Value = 3.5
$ws.Range("AU13").Value = 7
$ws.Range("AV13").Value = 41
$ws.Range("AW13").Value = 6
$ws.Range("AX13").Value = 19
$ws.Range("AY13").Value = 21
$ws.Range("AZ13").Value = 51
$ws.Range("BA13").Value = 67
$ws.Range("BB13").Value = 126
$ws.Range("BC13").Value = 351
$ws.Range("BD13").Value = 151

# --- Row 14: USA MLS FC Cincinnati vs New York City replaced by Seattle Sounders vs Houston Dynamo ---
$ws.Range("A14").Value = "CC5M2P9d"
$ws.Range("B14").Value = "28/10/2024"
$ws.Range("C14").Value = "21:50"
$ws.Range("D14").Value = "USA - MLS"
$ws.Range("E14").Value = "Seattle Sounders"
$ws.Range("F14").Value = "Houston Dynamo"
$ws.Range("G14").Value = 1.8
$ws.Range("H14").Value = 3.6
$ws.Range("I14").Value = 4.75
$ws.Range("J14").Value = 2.4
$ws.Range("K14").Value = 2.1
$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 9
$ws.Range("O14").Value = 1.33
$ws.Range("P14").Value = 3.25
$ws.Range("Q14").Value = 2.05
$ws.Range("R14").Value = 1.75
$ws.Range("S14").Value = 1.44
$ws.Range("T14").Value = 2.63
$ws.Range("U14").Value = 1.95
$ws.Range("V14").Value = 1.8
$ws.Range("W14").Value = 6.5
$ws.Range("X14").Value = 8
$ws.Range("Y14").Value = 8.5
$ws.Range("Z14").Value = 15
$ws.Range("AA14").Value = 15
$ws.Range("AB14").Value = 29
$ws.Range("AC14").Value = 9
$ws.Range("AD14").Value = 6.5
$ws.Range("AE14").Value = 17
$ws.Range("AF14").Value = 51
$ws.Range("AG14").Value = 11
$ws.Range("AH14").Value = 23
$ws.Range("AI14").Value = 15
$ws.Range("AJ14").Value = 51
$ws.Range("AK14").Value = 41
$ws.Range("AL14").Value = 41
$ws.Range("AM14").Value = 351
$ws.Range("AN14").Value = 3.75
$ws.Range("AO14").Value = 9.5
$ws.Range("AP14").Value = 21
$ws.Range("AQ14").Value = 34
$ws.Range("AR14").Value = 51
$ws.Range("AS14").Value = 151
$ws.Range("AT14").Value = 2.63
$ws.Range("AU14").Value = 8.5
$ws.Range("AV14").Value = 51
$ws.Range("AW14").Value = 6
$ws.Range("AX14").Value = 26
$ws.Range("AY14").Value = 34
$ws.Range("AZ14").Value = 81
$ws.Range("BA14").Value = 126
$ws.Range("BB14").Value = 251
$ws.Range("BC14").Value = 151
$ws.Range("BD14").Value = 151
